$wb = $excel.ActiveWorkbook

# --- Productdata sheet: scale safety-stock columns (D, F, I) by 0.0004 ---
$ws = $wb.Worksheets.Item("Productdata")

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 4).Value2 * 0.0004
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 6).Value2 * 0.0004
    $ws.Cells.Item($r, 9).Value2 = $ws.Cells.Item($r, 9).Value2 * 0.0004
}

# --- ForcastedStandardDeviation sheet: zero out B9:E11 ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws2.Range("B9:E11").Value2 = 0
